$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary field updates ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 10

# --- Insert a new line-item row before the current row 18 ---
# (old rows 18-25 shift down to 19-26; row 26 becomes the new TOTAL row)
$ws.Rows.Item(18).Insert()

# Row-insert shifts cell formatting down together with the row above it,
# which leaves the even/odd banding (style ids 9/10/11 vs 12/13/14) out of
# sync with absolute row position. Re-stamp the banding so it keeps
# alternating purely by row number, exactly like rows 16 (A-style) and 17
# (B-style) already do, untouched by the insert.
$ws.Range("A16:H16").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)
$ws.Range("A20:H20").PasteSpecial(-4122)
$ws.Range("A22:H22").PasteSpecial(-4122)
$ws.Range("A24:H24").PasteSpecial(-4122)

$ws.Range("A17:H17").Copy()
$ws.Range("A19:H19").PasteSpecial(-4122)
$ws.Range("A21:H21").PasteSpecial(-4122)
$ws.Range("A23:H23").PasteSpecial(-4122)
$ws.Range("A25:H25").PasteSpecial(-4122)

# --- Row 16 (unchanged content; pricing zeroed out) ---
$ws.Range("H16").Value = 0

# --- Row 17 (unchanged content; pricing zeroed out) ---
$ws.Range("H17").Value = 0

# --- Row 18: new SVW line item ---
$ws.Range("A18").Value = "Point 01"
$ws.Range("B18").Value = "SVW-2-TP-CUS-CC"
$ws.Range("C18").Value = "Inst"
$ws.Range("D18").Value = "SVW,#2 AWG,Trip,Copper Str,Corr/Comm"
$ws.Range("E18").Value = "FT"
$ws.Range("F18").Value = 15
$ws.Range("H18").Value = 0

# --- Row 19: previously row 18 (CON-2-AAI-3-P) ---
$ws.Range("A19").Value = "Point 01"
$ws.Range("B19").Value = "CON-2-AAI-3-P"
$ws.Range("C19").Value = "Inst"
$ws.Range("D19").Value = "CON,#2 AWG,AA Corr,Three,Poly"
$ws.Range("E19").Value = "FT"
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = 0

# --- Row 20: previously row 19 (CNC-NTI-10) ---
$ws.Range("A20").Value = "Point 01"
$ws.Range("B20").Value = "CNC-NTI-10"
$ws.Range("C20").Value = "Inst"
$ws.Range("D20").Value = "CNC,splice Non-Tension Insul,336-1033"
$ws.Range("E20").Value = "EA"
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 0

# --- Row 21: previously row 20 (Point 09, PLA-DLOC) ---
$ws.Range("A21").Value = "Point 09"
$ws.Range("B21").Value = "PLA-DLOC"
$ws.Range("C21").Value = "Inst"
$ws.Range("D21").Value = "PLA,Difficult Location"
$ws.Range("E21").Value = "EA"
$ws.Range("F21").Value = 6
$ws.Range("H21").Value = 0

# --- Row 22: previously row 21 (Point 04, PLA-DLOC) ---
$ws.Range("A22").Value = "Point 04"
$ws.Range("B22").Value = "PLA-DLOC"
$ws.Range("C22").Value = "Inst"
$ws.Range("D22").Value = "PLA,Difficult Location"
$ws.Range("E22").Value = "EA"
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 0

# --- Row 23: previously row 22 (Point 05, PLA-DLOC) ---
$ws.Range("A23").Value = "Point 05"
$ws.Range("B23").Value = "PLA-DLOC"
$ws.Range("C23").Value = "Inst"
$ws.Range("D23").Value = "PLA,Difficult Location"
$ws.Range("E23").Value = "EA"
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 0

# --- Row 24: previously row 23 (Point 06, CNC-HTA-40) ---
$ws.Range("A24").Value = "Point 06"
$ws.Range("B24").Value = "CNC-HTA-40"
$ws.Range("C24").Value = "Inst"
$ws.Range("D24").Value = "Compression Connector H-Tap Assembly 4/0"
$ws.Range("E24").Value = "EA"
$ws.Range("F24").Value = 21
$ws.Range("H24").Value = 0

# --- Row 25: previously row 24 (Point 08, PLA-DLOC) ---
$ws.Range("A25").Value = "Point 08"
$ws.Range("B25").Value = "PLA-DLOC"
$ws.Range("C25").Value = "Inst"
$ws.Range("D25").Value = "PLA,Difficult Location"
$ws.Range("E25").Value = "EA"
$ws.Range("F25").Value = 6
$ws.Range("H25").Value = 0

# --- Total row (now row 26 after the insert) ---
$ws.Range("H26").Value = 0
